$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Have you ever been infected by COVID-19?'
$ws.Range("C4").Value = 'Have any other family members been infected by COVID-19?'
$ws.Range("C6").Value = 'Due to COVID-19, did your family have to migrate to a different location?'
$ws.Range("C8").Value = 'Did you experience any of the following due to the economic downturn brought about by covid-19?'
$ws.Range("C9").Value = 'Due to have to any of the following due to COVID-19?'
$ws.Range("C10").Value = 'When was the last time you received your full salary (or full income)?'
$ws.Range("C16").Value = 'How many family members, including yourself, depended on your income before Covid-19?'
$ws.Range("C17").Value = 'How many family members, including yourself, depend on your income after Covid-19?'
$ws.Range("C18").Value = 'How many of your family members’ basic needs can be met by your current income?'
$ws.Range("C19").Value = 'Compared to the end of 2019, how much savings do you have now?'
$ws.Range("C20").Value = 'Which of the following assets and amenities did your family have before Covid-19?'
$ws.Range("C21").Value = 'Which of the following assets and amenities does your family have after Covid-19?'
$ws.Range("C22").Value = 'Do you have an outstanding debt?'
$ws.Range("C23").Value = 'Did you learn any new skills when the your tourism-related job/profession was halted due to Covid-19?'
$ws.Range("C26").Value = 'Who provided this training?'
$ws.Range("C27").Value = 'Which of the following sources do you mainly use to get information about Covid-19? (Choose three main sources)'
$ws.Range("C28").Value = 'Which of the following would help in reducing the effects of COVID-19 on you?'
$ws.Range("C29").Value = 'How much longer will COVID-19 affect your livelihood? Provide your best estimate.'
$ws.Range("C30").Value = 'How long will it take for your employment to return to normal after the tourism sector is up and running again? Provide your best estimate.'
$ws.Range("C31").Value = 'If the current situation continues for the next 6 months, which of the following difficulties will you face?'
$ws.Range("C32").Value = 'Do you have the necessary collateral/securities to secure additional loans?'
$ws.Range("C33").Value = 'What are the major challenges for the revival of the tourism sub-sector you are involved in?'
